$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("54").Insert()

$ws.Range("A54").Value = 7
$ws.Range("B54").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C54").Value = "Ñuble"
$ws.Range("D54").Value = 44966
$ws.Range("E54").Value = 16
$ws.Range("F54").Value = 100112030
$ws.Range("G54").Value = "Poroto granado"
$ws.Range("H54").Value = "Sin especificar"
$ws.Range("I54").Value = "Primera"
$ws.Range("J54").Value = 50
$ws.Range("K54").Value = 30000
$ws.Range("L54").Value = 30000
$ws.Range("M54").Value = 30000
$ws.Range("N54").Value = "$/saco 25 kilos"
$ws.Range("O54").Value = "Región del Maule"
$ws.Range("P54").Value = 1200
$ws.Range("Q54").Value = 25
$ws.Range("R54").Value = "Hortaliza"
